# Adicionado filtro de dados
# Adds a new user row (Charles / funcionario) to the usuarios sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Charles"
$ws.Range("B4").Value = "`$2b`$12`$QQjbf5ap2t8m6gHqIJ5WxeUHn4VJFcDblqQShQWC92oupeMbk61S."
$ws.Range("C4").Value = "funcionario"
